$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "DaysInCycle" in J1, matching style of existing headers (A1:I1)
$ws.Range("A1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "DaysInCycle"

# Set column width for column J (target stored width = 14; COM ColumnWidth
# uses a slightly different unit, so compensate to land exactly on 14)
$ws.Columns.Item(10).ColumnWidth = 13.14

# Update A2 value from 180 to 100
$ws.Range("A2").Value = 100

# Add new value for DaysInCycle in J2
$ws.Range("J2").Value = 2

# Update selection to reflect new active cell
$ws.Range("J3").Select()
